$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update task descriptions that changed text (weeks 3-6, column C)
# Order matters for shared-string table layout, matching original authoring order
$ws.Range("C7").Value = "Shop filters - utánaolvasni, hogyan szokás elkészíteni, backend queryk? Megvalósítása, webes kliensben megvalósítás - shopban filter felület, filter service?"
$ws.Range("C6").Value = "Felhasználókezelés, kosár kezelése - felület létrehozás - bejelentkezés, profil, rendeléseim képernyő"
$ws.Range("C4").Value = "Üzleti folyamat feltérképezés, Felhasználókezelés, kosár kezelése, utánanézni, mik a lehetőségek"
$ws.Range("C5").Value = "Felhasználókezelés, kosár kezelése - backend megvalósítás, backend - NRT hozzáadása, Domain fejlesztése, refaktorálás"

# Widen column C to fit the longer text
$ws.Columns.Item(3).ColumnWidth = 138.7109375

# Update the active selection to C5
$ws.Range("C5").Select()
